$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare new row 16 by copying the format of row 15 (so styles match data rows) ---
$ws.Range("A15:AQ15").Copy()
$ws.Range("A16:AQ16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the red-font formatting previously applied to F2:G2 and F3:G3 ---
# (copy the plain left-aligned format used elsewhere in column F/G, e.g. F4, onto these cells)
$ws.Range("F4").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Write the refreshed report data (rows 2-16) ---
$row2 = New-Object "object[,]" 1,43
$row2[0,0] = 22443138
$row2[0,1] = 23286164
$row2[0,2] = "ACC"
$row2[0,3] = 45831.58079861111
$row2[0,4] = 45831.58729166666
$row2[0,5] = "0 a 02 dias"
$row2[0,6] = "0 a 02 dias"
$row2[0,7] = "QZUVKB"
$row2[0,8] = "MANUAL"
$row2[0,9] = "ALDA MARINA DE CAMPOS MELO"
$row2[0,10] = "Juliana Cardoso"
$row2[0,11] = "Juliana Cardoso"
$row2[0,12] = 45828.58055555556
$row2[0,13] = "-"
$row2[0,14] = "-"
$row2[0,15] = "OFF LINE"
$row2[0,16] = "Cartão de crédito"
$row2[0,17] = "Cartão de crédito"
$row2[0,18] = "Aéreo"
$row2[0,19] = "N"
$row2[0,20] = "Grupo Kontrip"
$row2[0,21] = "Vpr Consultoria, Eventos, Servicos e Treinamentos Ltda"
$row2[0,22] = "-"
$row2[0,23] = "Gol Linhas Aereas"
$row2[0,24] = 2139988704
$row2[0,25] = "KONTRIP"
$row2[0,26] = "-"
$row2[0,27] = 2015.77
$row2[0,28] = 59.95
$row2[0,29] = 0
$row2[0,30] = 0
$row2[0,31] = 0
$row2[0,32] = 0
$row2[0,33] = 120.95
$row2[0,34] = "-"
$row2[0,35] = "35Ocorreu a seguinte exceção ao gerar a ordem de venda: Ocorreu a seguinte exceção ao inserir o item da ordem de venda:  Faltou informar rateio de "
$row2[0,36] = "KONTRIP"
$row2[0,37] = "Falta de informação Gerencial"
$row2[0,38] = "Rateio de centro de custo/projeto"
$row2[0,39] = "Dados Gerenciais"
$row2[0,40] = "Qualidade dos dados"
$row2[0,41] = "KONTRIP VIAGENS"
$row2[0,42] = "Operações - KONTRIP"
$ws.Range("A2:AQ2").Value = $row2

$row3 = New-Object "object[,]" 1,43
$row3[0,0] = 22443418
$row3[0,1] = 23286368
$row3[0,2] = "ACC"
$row3[0,3] = 45831.61105324074
$row3[0,4] = 45831.61527777778
$row3[0,5] = "0 a 02 dias"
$row3[0,6] = "0 a 02 dias"
$row3[0,7] = "HRTRXM"
$row3[0,8] = "MANUAL"
$row3[0,9] = "SILVIO CRESPO"
$row3[0,10] = "Juliana Cardoso"
$row3[0,11] = "Juliana Cardoso"
$row3[0,12] = 45828.61041666667
$row3[0,13] = "-"
$row3[0,14] = "-"
$row3[0,15] = "OFF LINE"
$row3[0,16] = "Cartão de crédito"
$row3[0,17] = "Cartão de crédito"
$row3[0,18] = "Aéreo"
$row3[0,19] = "N"
$row3[0,20] = "Grupo Kontrip"
$row3[0,21] = "Grana Capital Tecnologia Em Investimentos Ltda"
$row3[0,22] = "-"
$row3[0,23] = "Gol Linhas Aereas"
$row3[0,24] = 2139949448
$row3[0,25] = "KONTRIP"
$row3[0,26] = "-"
$row3[0,27] = 198.7
$row3[0,28] = 34.11
$row3[0,29] = 0
$row3[0,30] = 0
$row3[0,31] = 0
$row3[0,32] = 0
$row3[0,33] = 11.92
$row3[0,34] = "-"
$row3[0,35] = "35Ocorreu a seguinte exceção ao gerar a ordem de venda: Ocorreu a seguinte exceção ao inserir o item da ordem de venda:  Faltou informar rateio de "
$row3[0,36] = "KONTRIP"
$row3[0,37] = "Falta de informação Gerencial"
$row3[0,38] = "Rateio de centro de custo/projeto"
$row3[0,39] = "Dados Gerenciais"
$row3[0,40] = "Qualidade dos dados"
$row3[0,41] = "KONTRIP VIAGENS"
$row3[0,42] = "Operações - KONTRIP"
$ws.Range("A3:AQ3").Value = $row3

$row4 = New-Object "object[,]" 1,43
$row4[0,0] = 22443698
$row4[0,1] = 23286625
$row4[0,2] = "ACC"
$row4[0,3] = 45831.63736111111
$row4[0,4] = 45831.63805555556
$row4[0,5] = "0 a 02 dias"
$row4[0,6] = "0 a 02 dias"
$row4[0,7] = "UGDCFT"
$row4[0,8] = "EBOOKING"
$row4[0,9] = "FERNANDES/DANIEL"
$row4[0,10] = "Kontrip"
$row4[0,11] = "Kontrip"
$row4[0,12] = 45831.59166666667
$row4[0,13] = "-"
$row4[0,14] = "-"
$row4[0,15] = "OFF LINE"
$row4[0,16] = "Cartão de crédito"
$row4[0,17] = "Cartão de crédito"
$row4[0,18] = "Aéreo"
$row4[0,19] = "N"
$row4[0,20] = "Independente"
$row4[0,21] = "Tap Air Portugal"
$row4[0,22] = "-"
$row4[0,23] = "Gol Linhas Aereas"
$row4[0,24] = 2140079741
$row4[0,25] = "KONTRIP"
$row4[0,26] = "-"
$row4[0,27] = 928.1
$row4[0,28] = 50.92
$row4[0,29] = 0
$row4[0,30] = 0
$row4[0,31] = 0
$row4[0,32] = 0
$row4[0,33] = 55.69
$row4[0,34] = "obs"
$row4[0,35] = "Centro de custo não preenchido! (ACC01)"
$row4[0,36] = "KONTRIP"
$row4[0,37] = "Centro de custo"
$row4[0,38] = "Falta de informação Gerencial"
$row4[0,39] = "Dados do Fornecedor"
$row4[0,40] = "Qualidade dos dados"
$row4[0,41] = "KONTRIP VIAGENS"
$row4[0,42] = "Operações - KONTRIP"
$ws.Range("A4:AQ4").Value = $row4

$row5 = New-Object "object[,]" 1,43
$row5[0,0] = 22443998
$row5[0,1] = 23286896
$row5[0,2] = "ACC01"
$row5[0,3] = 45831.67209490741
$row5[0,4] = 45831.67363425926
$row5[0,5] = "0 a 02 dias"
$row5[0,6] = "0 a 02 dias"
$row5[0,7] = "LA9574694UPEM"
$row5[0,8] = "EBOOKING"
$row5[0,9] = "rui pinhal/mario"
$row5[0,10] = "Kontrip"
$row5[0,11] = "Kontrip"
$row5[0,12] = 45831.59652777778
$row5[0,13] = "-"
$row5[0,14] = "-"
$row5[0,15] = "OFF LINE"
$row5[0,16] = "Cartão de crédito"
$row5[0,17] = "Cartão de crédito"
$row5[0,18] = "Aéreo"
$row5[0,19] = "N"
$row5[0,20] = "Independente"
$row5[0,21] = "Tap Air Portugal"
$row5[0,22] = "-"
$row5[0,23] = "Latam Airlines Brasil"
$row5[0,24] = 2238221759
$row5[0,25] = "KONTRIP"
$row5[0,26] = "-"
$row5[0,27] = 919.11
$row5[0,28] = 82.36
$row5[0,29] = 0
$row5[0,30] = 0
$row5[0,31] = 0
$row5[0,32] = 0
$row5[0,33] = 13.79
$row5[0,34] = "obs"
$row5[0,35] = "Centro de custo não preenchido! (ACC01)"
$row5[0,36] = "KONTRIP"
$row5[0,37] = "Centro de custo"
$row5[0,38] = "Falta de informação Gerencial"
$row5[0,39] = "Dados do Fornecedor"
$row5[0,40] = "Qualidade dos dados"
$row5[0,41] = "KONTRIP VIAGENS"
$row5[0,42] = "Operações - KONTRIP"
$ws.Range("A5:AQ5").Value = $row5

$row6 = New-Object "object[,]" 1,43
$row6[0,0] = 22442897
$row6[0,1] = 23285884
$row6[0,2] = "ACC"
$row6[0,3] = 45831.55122685185
$row6[0,4] = 45831.55306712963
$row6[0,5] = "0 a 02 dias"
$row6[0,6] = "0 a 02 dias"
$row6[0,7] = "FDUEIN"
$row6[0,8] = "EBOOKING"
$row6[0,9] = "Ribeiro Guth/Alberto"
$row6[0,10] = "Kontrip"
$row6[0,11] = "Kontrip"
$row6[0,12] = 45831.54930555556
$row6[0,13] = "-"
$row6[0,14] = "-"
$row6[0,15] = "OFF LINE"
$row6[0,16] = "Cartão de crédito"
$row6[0,17] = "Cartão de crédito"
$row6[0,18] = "Aéreo"
$row6[0,19] = "N"
$row6[0,20] = "Grupo Kontrip"
$row6[0,21] = "Matterhorn Infraestrutura Gestao de Investimentos Ltda"
$row6[0,22] = "-"
$row6[0,23] = "Latam Airlines Brasil"
$row6[0,24] = 2238191518
$row6[0,25] = "KONTRIP"
$row6[0,26] = "-"
$row6[0,27] = 2498.21
$row6[0,28] = 289.52
$row6[0,29] = 0
$row6[0,30] = 0
$row6[0,31] = 0
$row6[0,32] = 0
$row6[0,33] = 40
$row6[0,34] = "obs"
$row6[0,35] = "Centro de custo não preenchido! (ACC01) Solicitante não preenchido! (ACC01)"
$row6[0,36] = "KONTRIP"
$row6[0,37] = "Mais de um campo não preenchido"
$row6[0,38] = "Falta de informação Gerencial"
$row6[0,39] = "Dados do Fornecedor"
$row6[0,40] = "Qualidade dos dados"
$row6[0,41] = "KONTRIP VIAGENS"
$row6[0,42] = "Operações - KONTRIP"
$ws.Range("A6:AQ6").Value = $row6

$row7 = New-Object "object[,]" 1,43
$row7[0,0] = 22443137
$row7[0,1] = 23286108
$row7[0,2] = "ACC"
$row7[0,3] = 45831.58091435185
$row7[0,4] = 45831.58263888889
$row7[0,5] = "0 a 02 dias"
$row7[0,6] = "0 a 02 dias"
$row7[0,7] = "LA9577948KSFD"
$row7[0,8] = "EBOOKING"
$row7[0,9] = "ARAUJO/BENITO"
$row7[0,10] = "Kontrip"
$row7[0,11] = "Kontrip"
$row7[0,12] = 45831.43680555555
$row7[0,13] = "-"
$row7[0,14] = "-"
$row7[0,15] = "OFF LINE"
$row7[0,16] = "Cartão de crédito"
$row7[0,17] = "Cartão de crédito"
$row7[0,18] = "Aéreo"
$row7[0,19] = "N"
$row7[0,20] = "Grupo Kontrip"
$row7[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row7[0,22] = "-"
$row7[0,23] = "Latam Airlines Brasil"
$row7[0,24] = 2238198211
$row7[0,25] = "KONTRIP"
$row7[0,26] = "-"
$row7[0,27] = 1161.19
$row7[0,28] = 48.16
$row7[0,29] = 0
$row7[0,30] = 0
$row7[0,31] = 0
$row7[0,32] = 0
$row7[0,33] = 487.71
$row7[0,34] = "obs"
$row7[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row7[0,36] = "KONTRIP"
$row7[0,37] = "Bilhete duplicado"
$row7[0,38] = "Bilhete Já Contabilizado"
$row7[0,39] = "Duplicidade de Contabilização"
$row7[0,40] = "Qualidade dos dados"
$row7[0,41] = "KONTRIP VIAGENS"
$row7[0,42] = "Conciliação aérea"
$ws.Range("A7:AQ7").Value = $row7

$row8 = New-Object "object[,]" 1,43
$row8[0,0] = 22443137
$row8[0,1] = 23286109
$row8[0,2] = "ACC02"
$row8[0,3] = 45831.58091435185
$row8[0,4] = 45831.58263888889
$row8[0,5] = "0 a 02 dias"
$row8[0,6] = "0 a 02 dias"
$row8[0,7] = "LA9577948KSFD"
$row8[0,8] = "EBOOKING"
$row8[0,9] = "SILVA/SUELLEN"
$row8[0,10] = "Kontrip"
$row8[0,11] = "Kontrip"
$row8[0,12] = 45831.43680555555
$row8[0,13] = "-"
$row8[0,14] = "-"
$row8[0,15] = "OFF LINE"
$row8[0,16] = "Cartão de crédito"
$row8[0,17] = "Cartão de crédito"
$row8[0,18] = "Aéreo"
$row8[0,19] = "N"
$row8[0,20] = "Grupo Kontrip"
$row8[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row8[0,22] = "-"
$row8[0,23] = "Latam Airlines Brasil"
$row8[0,24] = 2238198212
$row8[0,25] = "KONTRIP"
$row8[0,26] = "-"
$row8[0,27] = 1161.22
$row8[0,28] = 48.16
$row8[0,29] = 0
$row8[0,30] = 0
$row8[0,31] = 0
$row8[0,32] = 0
$row8[0,33] = 0
$row8[0,34] = "obs"
$row8[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row8[0,36] = "KONTRIP"
$row8[0,37] = "Bilhete duplicado"
$row8[0,38] = "Bilhete Já Contabilizado"
$row8[0,39] = "Duplicidade de Contabilização"
$row8[0,40] = "Qualidade dos dados"
$row8[0,41] = "KONTRIP VIAGENS"
$row8[0,42] = "Conciliação aérea"
$ws.Range("A8:AQ8").Value = $row8

$row9 = New-Object "object[,]" 1,43
$row9[0,0] = 22443137
$row9[0,1] = 23286110
$row9[0,2] = "ACC03"
$row9[0,3] = 45831.58091435185
$row9[0,4] = 45831.58263888889
$row9[0,5] = "0 a 02 dias"
$row9[0,6] = "0 a 02 dias"
$row9[0,7] = "LA9577948KSFD"
$row9[0,8] = "EBOOKING"
$row9[0,9] = "Marques/Jonnata"
$row9[0,10] = "Kontrip"
$row9[0,11] = "Kontrip"
$row9[0,12] = 45831.43680555555
$row9[0,13] = "-"
$row9[0,14] = "-"
$row9[0,15] = "OFF LINE"
$row9[0,16] = "Cartão de crédito"
$row9[0,17] = "Cartão de crédito"
$row9[0,18] = "Aéreo"
$row9[0,19] = "N"
$row9[0,20] = "Grupo Kontrip"
$row9[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row9[0,22] = "-"
$row9[0,23] = "Latam Airlines Brasil"
$row9[0,24] = 2238198213
$row9[0,25] = "KONTRIP"
$row9[0,26] = "-"
$row9[0,27] = 1161.22
$row9[0,28] = 48.16
$row9[0,29] = 0
$row9[0,30] = 0
$row9[0,31] = 0
$row9[0,32] = 0
$row9[0,33] = 0
$row9[0,34] = "obs"
$row9[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row9[0,36] = "KONTRIP"
$row9[0,37] = "Bilhete duplicado"
$row9[0,38] = "Bilhete Já Contabilizado"
$row9[0,39] = "Duplicidade de Contabilização"
$row9[0,40] = "Qualidade dos dados"
$row9[0,41] = "KONTRIP VIAGENS"
$row9[0,42] = "Conciliação aérea"
$ws.Range("A9:AQ9").Value = $row9

$row10 = New-Object "object[,]" 1,43
$row10[0,0] = 22443137
$row10[0,1] = 23286111
$row10[0,2] = "ACC04"
$row10[0,3] = 45831.58091435185
$row10[0,4] = 45831.58263888889
$row10[0,5] = "0 a 02 dias"
$row10[0,6] = "0 a 02 dias"
$row10[0,7] = "LA9577948KSFD"
$row10[0,8] = "EBOOKING"
$row10[0,9] = "JANUARIO/BRUNO"
$row10[0,10] = "Kontrip"
$row10[0,11] = "Kontrip"
$row10[0,12] = 45831.43680555555
$row10[0,13] = "-"
$row10[0,14] = "-"
$row10[0,15] = "OFF LINE"
$row10[0,16] = "Cartão de crédito"
$row10[0,17] = "Cartão de crédito"
$row10[0,18] = "Aéreo"
$row10[0,19] = "N"
$row10[0,20] = "Grupo Kontrip"
$row10[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row10[0,22] = "-"
$row10[0,23] = "Latam Airlines Brasil"
$row10[0,24] = 2238198214
$row10[0,25] = "KONTRIP"
$row10[0,26] = "-"
$row10[0,27] = 1161.22
$row10[0,28] = 48.16
$row10[0,29] = 0
$row10[0,30] = 0
$row10[0,31] = 0
$row10[0,32] = 0
$row10[0,33] = 0
$row10[0,34] = "obs"
$row10[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row10[0,36] = "KONTRIP"
$row10[0,37] = "Bilhete duplicado"
$row10[0,38] = "Bilhete Já Contabilizado"
$row10[0,39] = "Duplicidade de Contabilização"
$row10[0,40] = "Qualidade dos dados"
$row10[0,41] = "KONTRIP VIAGENS"
$row10[0,42] = "Conciliação aérea"
$ws.Range("A10:AQ10").Value = $row10

$row11 = New-Object "object[,]" 1,43
$row11[0,0] = 22443137
$row11[0,1] = 23286112
$row11[0,2] = "ACC05"
$row11[0,3] = 45831.58091435185
$row11[0,4] = 45831.58263888889
$row11[0,5] = "0 a 02 dias"
$row11[0,6] = "0 a 02 dias"
$row11[0,7] = "LA9577948KSFD"
$row11[0,8] = "EBOOKING"
$row11[0,9] = "ALBINO/GUSTAVO"
$row11[0,10] = "Kontrip"
$row11[0,11] = "Kontrip"
$row11[0,12] = 45831.43680555555
$row11[0,13] = "-"
$row11[0,14] = "-"
$row11[0,15] = "OFF LINE"
$row11[0,16] = "Cartão de crédito"
$row11[0,17] = "Cartão de crédito"
$row11[0,18] = "Aéreo"
$row11[0,19] = "N"
$row11[0,20] = "Grupo Kontrip"
$row11[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row11[0,22] = "-"
$row11[0,23] = "Latam Airlines Brasil"
$row11[0,24] = 2238198215
$row11[0,25] = "KONTRIP"
$row11[0,26] = "-"
$row11[0,27] = 1161.22
$row11[0,28] = 48.16
$row11[0,29] = 0
$row11[0,30] = 0
$row11[0,31] = 0
$row11[0,32] = 0
$row11[0,33] = 0
$row11[0,34] = "obs"
$row11[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row11[0,36] = "KONTRIP"
$row11[0,37] = "Bilhete duplicado"
$row11[0,38] = "Bilhete Já Contabilizado"
$row11[0,39] = "Duplicidade de Contabilização"
$row11[0,40] = "Qualidade dos dados"
$row11[0,41] = "KONTRIP VIAGENS"
$row11[0,42] = "Conciliação aérea"
$ws.Range("A11:AQ11").Value = $row11

$row12 = New-Object "object[,]" 1,43
$row12[0,0] = 22443137
$row12[0,1] = 23286113
$row12[0,2] = "ACC06"
$row12[0,3] = 45831.58091435185
$row12[0,4] = 45831.58263888889
$row12[0,5] = "0 a 02 dias"
$row12[0,6] = "0 a 02 dias"
$row12[0,7] = "LA9577948KSFD"
$row12[0,8] = "EBOOKING"
$row12[0,9] = "CAMARGO/JULIANA"
$row12[0,10] = "Kontrip"
$row12[0,11] = "Kontrip"
$row12[0,12] = 45831.43680555555
$row12[0,13] = "-"
$row12[0,14] = "-"
$row12[0,15] = "OFF LINE"
$row12[0,16] = "Cartão de crédito"
$row12[0,17] = "Cartão de crédito"
$row12[0,18] = "Aéreo"
$row12[0,19] = "N"
$row12[0,20] = "Grupo Kontrip"
$row12[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row12[0,22] = "-"
$row12[0,23] = "Latam Airlines Brasil"
$row12[0,24] = 2238198216
$row12[0,25] = "KONTRIP"
$row12[0,26] = "-"
$row12[0,27] = 1161.22
$row12[0,28] = 48.16
$row12[0,29] = 0
$row12[0,30] = 0
$row12[0,31] = 0
$row12[0,32] = 0
$row12[0,33] = 0
$row12[0,34] = "obs"
$row12[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row12[0,36] = "KONTRIP"
$row12[0,37] = "Bilhete duplicado"
$row12[0,38] = "Bilhete Já Contabilizado"
$row12[0,39] = "Duplicidade de Contabilização"
$row12[0,40] = "Qualidade dos dados"
$row12[0,41] = "KONTRIP VIAGENS"
$row12[0,42] = "Conciliação aérea"
$ws.Range("A12:AQ12").Value = $row12

$row13 = New-Object "object[,]" 1,43
$row13[0,0] = 22443137
$row13[0,1] = 23286114
$row13[0,2] = "ACC07"
$row13[0,3] = 45831.58091435185
$row13[0,4] = 45831.58263888889
$row13[0,5] = "0 a 02 dias"
$row13[0,6] = "0 a 02 dias"
$row13[0,7] = "LA9577948KSFD"
$row13[0,8] = "EBOOKING"
$row13[0,9] = "RACHINSKI/RAFAEL"
$row13[0,10] = "Kontrip"
$row13[0,11] = "Kontrip"
$row13[0,12] = 45831.43680555555
$row13[0,13] = "-"
$row13[0,14] = "-"
$row13[0,15] = "OFF LINE"
$row13[0,16] = "Cartão de crédito"
$row13[0,17] = "Cartão de crédito"
$row13[0,18] = "Aéreo"
$row13[0,19] = "N"
$row13[0,20] = "Grupo Kontrip"
$row13[0,21] = "Associacao Colo de Deus e Santissima Virgem (a.c.d.s.v)"
$row13[0,22] = "-"
$row13[0,23] = "Latam Airlines Brasil"
$row13[0,24] = 2238198216
$row13[0,25] = "KONTRIP"
$row13[0,26] = "-"
$row13[0,27] = 1161.22
$row13[0,28] = 48.16
$row13[0,29] = 0
$row13[0,30] = 0
$row13[0,31] = 0
$row13[0,32] = 0
$row13[0,33] = 0
$row13[0,34] = "obs"
$row13[0,35] = "Verificação de bilhetes: Bilhete 2238198216 já sendo utilizado para este fornecedor."
$row13[0,36] = "KONTRIP"
$row13[0,37] = "Bilhete duplicado"
$row13[0,38] = "Bilhete Já Contabilizado"
$row13[0,39] = "Duplicidade de Contabilização"
$row13[0,40] = "Qualidade dos dados"
$row13[0,41] = "KONTRIP VIAGENS"
$row13[0,42] = "Conciliação aérea"
$ws.Range("A13:AQ13").Value = $row13

$row14 = New-Object "object[,]" 1,43
$row14[0,0] = 22444419
$row14[0,1] = 23287293
$row14[0,2] = "ACC01"
$row14[0,3] = 45831.72347222222
$row14[0,4] = 45831.72582175926
$row14[0,5] = "0 a 02 dias"
$row14[0,6] = "0 a 02 dias"
$row14[0,7] = "LA9571170HTTN"
$row14[0,8] = "EBOOKING"
$row14[0,9] = "Alves Magalhães Filho/Jorge"
$row14[0,10] = "Kontrip"
$row14[0,11] = "Kontrip"
$row14[0,12] = 45831.72222222222
$row14[0,13] = "-"
$row14[0,14] = "-"
$row14[0,15] = "OFF LINE"
$row14[0,16] = "Cartão de crédito"
$row14[0,17] = "Cartão de crédito"
$row14[0,18] = "Aéreo"
$row14[0,19] = "N"
$row14[0,20] = "Grupo Kontrip"
$row14[0,21] = "Lmaismaocubo Projetos e Consultoria Ltda"
$row14[0,22] = "-"
$row14[0,23] = "Latam Airlines Brasil"
$row14[0,24] = 2238234777
$row14[0,25] = "KONTRIP"
$row14[0,26] = "-"
$row14[0,27] = 419.02
$row14[0,28] = 100.18
$row14[0,29] = 0
$row14[0,30] = 0
$row14[0,31] = 0
$row14[0,32] = 0
$row14[0,33] = 40
$row14[0,34] = "obs"
$row14[0,35] = "Pnr já existente. A duplicidade de rloc é permitida apenas 6 meses após o último pnr emitido"
$row14[0,36] = "KONTRIP"
$row14[0,37] = "Duplicidade de RLOC"
$row14[0,38] = "Campo RLOC"
$row14[0,39] = "Duplicidade de Contabilização"
$row14[0,40] = "Qualidade dos dados"
$row14[0,41] = "KONTRIP VIAGENS"
$row14[0,42] = "Conciliação aérea"
$ws.Range("A14:AQ14").Value = $row14

$row15 = New-Object "object[,]" 1,43
$row15[0,0] = 22442491
$row15[0,1] = 23285512
$row15[0,2] = "ACC01"
$row15[0,3] = 45831.50450231481
$row15[0,4] = 45831.51125
$row15[0,5] = "0 a 02 dias"
$row15[0,6] = "0 a 02 dias"
$row15[0,7] = "OG584V"
$row15[0,8] = "EBOOKING"
$row15[0,9] = "DAVI ARAUJO"
$row15[0,10] = "Kontrip"
$row15[0,11] = "Kontrip"
$row15[0,12] = 45828.50416666667
$row15[0,13] = "-"
$row15[0,14] = "-"
$row15[0,15] = "OFF LINE"
$row15[0,16] = "Invoice"
$row15[0,17] = "Cartão convênio"
$row15[0,18] = "Aéreo"
$row15[0,19] = "N"
$row15[0,20] = "Grupo Kontrip"
$row15[0,21] = "Porto de Galinhas Convention e Visitors Bureau"
$row15[0,22] = "-"
$row15[0,23] = "Azul Linhas Aereas"
$row15[0,24] = 3023285511
$row15[0,25] = "KONTRIP"
$row15[0,26] = "-"
$row15[0,27] = 612.54
$row15[0,28] = 119.51
$row15[0,29] = 0
$row15[0,30] = 0
$row15[0,31] = 0
$row15[0,32] = 0
$row15[0,33] = 110.26
$row15[0,34] = "-"
$row15[0,35] = "36Faltou informar rateio de centro de custo/projeto abaixo da accounting"
$row15[0,36] = "KONTRIP"
$row15[0,37] = "Falta de informação Gerencial"
$row15[0,38] = "Rateio de centro de custo/projeto"
$row15[0,39] = "Dados Gerenciais"
$row15[0,40] = "Qualidade dos dados"
$row15[0,41] = "KONTRIP VIAGENS"
$row15[0,42] = "Operações - KONTRIP"
$ws.Range("A15:AQ15").Value = $row15

$row16 = New-Object "object[,]" 1,43
$row16[0,0] = 22442491
$row16[0,1] = 23285535
$row16[0,2] = "ACC02"
$row16[0,3] = 45831.50450231481
$row16[0,4] = 45831.51125
$row16[0,5] = "0 a 02 dias"
$row16[0,6] = "0 a 02 dias"
$row16[0,7] = "OG584V"
$row16[0,8] = "MANUAL"
$row16[0,9] = "DANIEL ARAUJO"
$row16[0,10] = "Juliana Cardoso"
$row16[0,11] = "Juliana Cardoso"
$row16[0,12] = 45828.50416666667
$row16[0,13] = "-"
$row16[0,14] = "-"
$row16[0,15] = "OFF LINE"
$row16[0,16] = "Invoice"
$row16[0,17] = "Cartão convênio"
$row16[0,18] = "Aéreo"
$row16[0,19] = "N"
$row16[0,20] = "Grupo Kontrip"
$row16[0,21] = "Porto de Galinhas Convention e Visitors Bureau"
$row16[0,22] = "-"
$row16[0,23] = "Azul Linhas Aereas"
$row16[0,24] = 3023285536
$row16[0,25] = "KONTRIP"
$row16[0,26] = "-"
$row16[0,27] = 612.54
$row16[0,28] = 119.51
$row16[0,29] = 0
$row16[0,30] = 0
$row16[0,31] = 0
$row16[0,32] = 0
$row16[0,33] = 0
$row16[0,34] = "-"
$row16[0,35] = "36Faltou informar rateio de centro de custo/projeto abaixo da accounting"
$row16[0,36] = "KONTRIP"
$row16[0,37] = "Falta de informação Gerencial"
$row16[0,38] = "Rateio de centro de custo/projeto"
$row16[0,39] = "Dados Gerenciais"
$row16[0,40] = "Qualidade dos dados"
$row16[0,41] = "KONTRIP VIAGENS"
$row16[0,42] = "Operações - KONTRIP"
$ws.Range("A16:AQ16").Value = $row16

